$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 183; existing rows 183:254 shift down to 184:255.
$ws.Rows("183:183").Insert()

# Populate the newly inserted row 183 with the new record.
$ws.Range("A183").Value = 7
$ws.Range("B183").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C183").Value = "Ñuble"
$ws.Range("D183").Value = Get-Date -Year 2022 -Month 11 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("E183").Value = 16
$ws.Range("F183").Value = 100112032
$ws.Range("G183").Value = "Zapallo italiano"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 160
$ws.Range("K183").Value = 8000
$ws.Range("L183").Value = 8500
$ws.Range("M183").Value = 8250
$ws.Range("N183").Value = "`$/caja 50 unidades"
$ws.Range("O183").Value = "Región de O'Higgins"
$ws.Range("P183").Value = 165
$ws.Range("Q183").Value = 50
$ws.Range("R183").Value = "Hortaliza"
